$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "42.935.60") that Excel would otherwise
# auto-parse as numbers; force text storage, then restore the default "Normal"
# style so the cells keep looking exactly as they did before (no style changes).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.935.60'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '2.361.74'
$ws.Range("E3").Value = '  +1.65%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '302.28'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").Value = '95.36'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -0.63%  '

$ws.Range("D9").Value = '0.485'
$ws.Range("E9").Value = '  -1.51%  '

$ws.Range("D10").Value = '33.99'
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("D12").Value = '0.0784'
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("D13").Value = '18.32'
$ws.Range("E13").Value = '  -3.44%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.729.74'
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '6.70'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '2.382.11'
$ws.Range("E16").Value = '  +2.43%  '

$ws.Range("D17").Value = '0.790'
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").Value = '42.898.40'
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("D19").Value = '11.88'
$ws.Range("E19").Value = '  -2.78%  '

$ws.Range("E20").Value = '  +1.64%  '

$ws.Range("E21").Value = '  -0.61%  '

$ws.Range("D22").Value = '68.01'
$ws.Range("E22").Value = '  +0.28%  '

$ws.Range("D23").Value = '234.99'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("E24").Value = '  -4.34%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").Value = '24.43'
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").Value = '  +0.65%  '

$ws.Range("E29").Value = '  +1.78%  '

$ws.Range("D30").Value = '31.87'
$ws.Range("E30").Value = '  -1.02%  '

$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").Value = '4.99'
$ws.Range("E32").Value = '  -0.22%  '

$ws.Range("D33").Value = '17.48'
$ws.Range("E33").Value = '  -1.65%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0717'
$ws.Range("E34").Value = '  +2.23%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '130.17'
$ws.Range("E35").Value = '  -11.18%  '

$ws.Range("E37").Value = '  +1.70%  '

$ws.Range("D38").Value = '4.32'
$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("D39").Value = '2.81'
$ws.Range("E39").Value = '  +2.47%  '

$ws.Range("E40").Value = '  -2.02%  '

$ws.Range("D41").Value = '0.107'
$ws.Range("E41").Value = '  -0.78%  '

$ws.Range("D42").Value = '21.24'
$ws.Range("E42").Value = '  -3.58%  '

$ws.Range("D43").Value = '1.928.78'
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").Value = '0.0278'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("E45").Value = '  +3.01%  '

$ws.Range("D46").Value = '9.17'
$ws.Range("E46").Value = '  -9.31%  '

$ws.Range("D47").Value = '2.70'
$ws.Range("E47").Value = '  -1.52%  '

$ws.Range("D48").Value = '2.591.73'
$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("E49").Value = '  +1.91%  '

$ws.Range("E50").Value = '  +1.51%  '

$ws.Range("D51").Value = '71.44'
$ws.Range("E51").Value = '  -1.32%  '

$ws.Range("D2:D51").Style = "Normal"
